# Updated cryptos list on Sat Mar 18 20:13:14 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# every coin row (rows 2-51) on the active worksheet with the latest
# scraped snapshot values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row new values. $null means that column did not change for that row.
$updates = @(
    @{ Row = 2; D = "27.407.55"; E = "  +2.44%  " },
    @{ Row = 3; D = "1.795.16"; E = "  +3.11%  " },
    @{ Row = 4; D = "1.004"; E = "  +0.34%  " },
    @{ Row = 5; D = "337.47"; E = "  +0.65%  " },
    @{ Row = 6; D = "1.001"; E = "  +0.23%  " },
    @{ Row = 7; D = "0.3798"; E = "  +1.40%  " },
    @{ Row = 8; D = $null; E = "  +1.71%  " },
    @{ Row = 9; D = "48.42"; E = "  -0.08%  " },
    @{ Row = 10; D = $null; E = "  +1.16%  " },
    @{ Row = 11; D = "0.07493"; E = "  +0.10%  " },
    @{ Row = 12; D = $null; E = "  +0.21%  " },
    @{ Row = 13; D = "22.07"; E = "  +7.81%  " },
    @{ Row = 14; D = "6.471"; E = "  +0.82%  " },
    @{ Row = 15; D = "1.797.82"; E = "  +3.19%  " },
    @{ Row = 16; D = "7.059"; E = "  -0.20%  " },
    @{ Row = 17; D = "0.00001100"; E = "  +1.90%  " },
    @{ Row = 18; D = "0.06645"; E = "  -1.26%  " },
    @{ Row = 19; D = "84.69"; E = "  +2.37%  " },
    @{ Row = 20; D = "1.000"; E = "  +0.25%  " },
    @{ Row = 21; D = "6.516"; E = "  +4.72%  " },
    @{ Row = 22; D = $null; E = "  +3.42%  " },
    @{ Row = 23; D = "27.389.99"; E = "  +2.38%  " },
    @{ Row = 24; D = "12.52"; E = "  -2.01%  " },
    @{ Row = 25; D = "2.433"; E = "  -0.77%  " },
    @{ Row = 26; D = "1.504"; E = "  +2.49%  " },
    @{ Row = 27; D = "2.559"; E = "  +5.09%  " },
    @{ Row = 28; D = "21.38"; E = "  +9.32%  " },
    @{ Row = 29; D = "151.91"; E = "  +0.14%  " },
    @{ Row = 30; D = "2.001.27"; E = "  +3.37%  " },
    @{ Row = 31; D = "134.05"; E = "  +1.12%  " },
    @{ Row = 32; D = "4.057"; E = "  -1.29%  " },
    @{ Row = 33; D = "6.119"; E = "  +0.90%  " },
    @{ Row = 34; D = "0.08684"; E = "  +0.28%  " },
    @{ Row = 35; D = $null; E = "  +2.55%  " },
    @{ Row = 36; D = "1.660"; E = "  -2.05%  " },
    @{ Row = 37; D = "5.445"; E = "  +0.08%  " },
    @{ Row = 38; D = "0.6899"; E = "  +9.96%  " },
    @{ Row = 39; D = "8.882"; E = $null },
    @{ Row = 40; D = "0.06371"; E = "  +1.49%  " },
    @{ Row = 41; D = $null; E = "  +1.35%  " },
    @{ Row = 42; D = $null; E = "  -0.99%  " },
    @{ Row = 43; D = "1.275"; E = "  +4.30%  " },
    @{ Row = 44; D = "14.48"; E = "  +0.58%  " },
    @{ Row = 45; D = "0.6438"; E = "  +5.47%  " },
    @{ Row = 46; D = "1.000"; E = "  +0.24%  " },
    @{ Row = 47; D = "3.868"; E = "  -1.50%  " },
    @{ Row = 48; D = "2.127"; E = "  +2.74%  " },
    @{ Row = 49; D = "130.08"; E = "  +0.58%  " },
    @{ Row = 50; D = "0.07193"; E = "  -0.35%  " },
    @{ Row = 51; D = "79.67"; E = "  +2.40%  " }
)

# Matches plain decimal-number-looking text (e.g. "1.004", "0.00001100",
# "337.47") as opposed to the "thousands.thousands.decimal"-style price
# strings (e.g. "27.407.55") that Excel can never mistake for a number.
# Cells whose new text could be auto-parsed as a number need to be marked
# as Text first so Excel keeps the exact original digits (incl. trailing
# zeros like "1.000") instead of silently re-typing them as a float.
$numberLike = '^[+-]?\d+(\.\d+)?$'

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        if ($u.D -match $numberLike) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
